$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F) column values for rows 2-6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1039
$ws1.Range("F3").Value = 186
$ws1.Range("F4").Value = 2379
$ws1.Range("F5").Value = 27
$ws1.Range("F6").Value = 521

# Sheet "全部类型" - update "想去人数" (F) column values for rows 4-8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1039
$ws4.Range("F5").Value = 186
$ws4.Range("F6").Value = 2379
$ws4.Range("F7").Value = 27
$ws4.Range("F8").Value = 521
